# Fixed Bento 80 Test scripts
# Appends/adjusts the trailing "order by ... LIMIT 100" clauses in the
# Neo4j query text stored in B2:B4 of the "startup" sheet, then updates
# the row heights (which Excel recalculates for the now-taller wrapped
# text) and moves the active selection to B4 to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2 (CasesTab query): add an ORDER BY / LIMIT clause on a new line ---
$b2 = $ws.Range("B2").Value2
$b2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Range("B2").Value = $b2

# --- B3 (SamplesTab query): add an ORDER BY / LIMIT clause on a new line ---
$b3 = $ws.Range("B3").Value2
$b3 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $b3

# --- B4 (FilesTab query): rewrite the existing "order by" tail ---
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "\n    order by f\.file_name$", "`n     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value = $b4

# The extra wrapped line pushes rows 2 and 3 taller (row 4 was already at
# Excel's row-height ceiling of 409.6 and stays there).
$ws.Rows(2).RowHeight = 360
$ws.Rows(3).RowHeight = 374.4

# Move / record the active selection on B4, matching the saved workbook.
$ws.Range("B4").Select() | Out-Null
